# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values for the
# 6f13978f-22ec-46d7-9ec0-473a44c01773 row on both the "zh-cn" and
# "de-de" worksheets. Row 5 (9b8aaa83-...) happens to share the very
# same timestamp text as row 3, so it must be updated too -- that is
# how the shared string table keeps a single de-duplicated entry for
# both rows, exactly like the original workbook.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-17 09:54:57"
$wsZhCn.Range("H3").Value = "2016-03-17 09:55:15"
$wsZhCn.Range("E5").Value = "2016-03-17 09:54:57"
$wsZhCn.Range("H5").Value = "2016-03-17 09:55:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-17 09:55:00"
$wsDeDe.Range("H3").Value = "2016-03-17 09:55:21"
$wsDeDe.Range("E5").Value = "2016-03-17 09:55:00"
$wsDeDe.Range("H5").Value = "2016-03-17 09:55:21"
